$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "40.636.57"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "2.367.15"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'310.26"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").Value = "'87.39"
$ws.Range("E6").Value = "  -6.59%  "
$ws.Range("D7").Value = "'0.528"
$ws.Range("E7").Value = "  -4.89%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  -4.89%  "
$ws.Range("D10").Value = "'0.0837"
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("D11").Value = "'30.74"
$ws.Range("E11").Value = "  -7.79%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "2.733.60"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("D14").Value = "'6.54"
$ws.Range("E14").Value = "  -5.77%  "
$ws.Range("D15").Value = "'15.01"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "2.360.09"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").Value = "'0.759"
$ws.Range("E17").Value = "  -5.75%  "
$ws.Range("D18").Value = "40.516.96"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "0.0₃0908"
$ws.Range("E19").Value = "  -4.90%  "
$ws.Range("D20").Value = "'6.13"
$ws.Range("E20").Value = "  -5.74%  "
$ws.Range("D21").Value = "'68.79"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "'10.78"
$ws.Range("E22").Value = "  -5.30%  "
$ws.Range("D23").Value = "'233.19"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -7.96%  "
$ws.Range("D27").Value = "'23.72"
$ws.Range("E27").Value = "  -6.07%  "
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("D29").Value = "'9.33"
$ws.Range("E29").Value = "  -4.53%  "
$ws.Range("D30").Value = "'33.91"
$ws.Range("E30").Value = "  -8.71%  "
$ws.Range("D31").Value = "'153.05"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'5.21"
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("D34").Value = "'0.0728"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.76"
$ws.Range("E37").Value = "  -5.84%  "
$ws.Range("D38").Value = "'15.82"
$ws.Range("E38").Value = "  -9.30%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.0995"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").Value = "'1.71"
$ws.Range("E40").Value = "  -9.32%  "
$ws.Range("D41").Value = "'3.86"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").Value = "'2.40"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("D43").Value = "1.953.50"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("D45").Value = "'17.61"
$ws.Range("E45").Value = "  -8.65%  "
$ws.Range("D46").Value = "'9.52"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = "  -9.18%  "
$ws.Range("D48").Value = "2.604.86"
$ws.Range("E48").Value = "  -4.11%  "
$ws.Range("D49").Value = "'93.07"
$ws.Range("E49").Value = "  -5.42%  "
$ws.Range("D50").Value = "'72.44"
$ws.Range("E50").Value = "  -5.91%  "
$ws.Range("D51").Value = "'50.50"
$ws.Range("E51").Value = "  -4.10%  "
